$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Período 01/2025)
$ws.Range("B2").Value = 587582.36
$ws.Range("C2").Value = 19586.07866666666
$ws.Range("D2").Value = 31800.92

# Row 3 (Período 02/2025)
$ws.Range("B3").Value = 763210.77
$ws.Range("C3").Value = 27257.5275
$ws.Range("D3").Value = 57356.01

# Row 4 (Período 03/2025)
$ws.Range("B4").Value = 670620.61
$ws.Range("C4").Value = 21632.9229032258
$ws.Range("D4").Value = 51663.12
$ws.Range("E4").Value = 7

# Row 5 (Período 04/2025)
$ws.Range("B5").Value = 511614.46
$ws.Range("C5").Value = 18271.945
$ws.Range("D5").Value = 38515.34
$ws.Range("E5").Value = 8
